$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.195.53"
$ws.Range("E2").Value = "  +0.75%  "
$ws.Range("D3").Value = "'3.504.01"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'604.00"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("D6").Value = "'172.70"
$ws.Range("E6").Value = "  +1.06%  "
$ws.Range("D7").Value = "'0.610"
$ws.Range("E7").Value = "  -0.58%  "
$ws.Range("D8").Value = "'3.501.43"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").Value = "'0.195"
$ws.Range("E10").Value = "  -1.98%  "
$ws.Range("D11").Value = "'7.20"
$ws.Range("E11").Value = "  +7.69%  "
$ws.Range("D12").Value = "'0.585"
$ws.Range("E12").Value = "  +0.81%  "
$ws.Range("D13").Value = "'46.02"
$ws.Range("E13").Value = "  -2.18%  "
$ws.Range("E14").Value = "  -0.87%  "
$ws.Range("D15").Value = "'4.075.89"
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").Value = "'613.76"
$ws.Range("D17").Value = "'8.26"
$ws.Range("E17").Value = "  -1.36%  "
$ws.Range("D18").Value = "'3.512.68"
$ws.Range("E18").Value = "  +0.25%  "
$ws.Range("D19").Value = "'70.228.45"
$ws.Range("E19").Value = "  +0.70%  "
$ws.Range("E20").Value = "  +0.94%  "
$ws.Range("D21").Value = "'17.49"
$ws.Range("E21").Value = "  +1.16%  "
$ws.Range("D22").Value = "'0.877"
$ws.Range("E22").Value = "  -0.48%  "
$ws.Range("D23").Value = "'9.07"
$ws.Range("E23").Value = "  -8.10%  "
$ws.Range("D24").Value = "'99.10"
$ws.Range("E24").Value = "  +3.34%  "
$ws.Range("D25").Value = "'15.55"
$ws.Range("E25").Value = "  -1.29%  "
$ws.Range("D26").Value = "'3.71"
$ws.Range("E26").Value = "  -3.18%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").Value = "'2.55"
$ws.Range("E28").Value = "  -1.24%  "
$ws.Range("D29").Value = "'34.03"
$ws.Range("E29").Value = "  +2.99%  "
$ws.Range("D30").Value = "'8.98"
$ws.Range("E30").Value = "  -1.82%  "
$ws.Range("B31").Value = "Stacks"
$ws.Range("C31").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D31").Value = "'2.97"
$ws.Range("E31").Value = "  -2.53%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'8.03"
$ws.Range("E32").Value = "  -4.06%  "
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").Value = "'651.92"
$ws.Range("E33").Value = "  +15.58%  "
$ws.Range("E34").Value = "  -4.40%  "
$ws.Range("D35").Value = "'6.79"
$ws.Range("E35").Value = "  -2.01%  "
$ws.Range("D36").Value = "'0.0993"
$ws.Range("E36").Value = "  -1.51%  "
$ws.Range("D37").Value = "'10.72"
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").Value = "'3.50"
$ws.Range("E38").Value = "  +0.92%  "
$ws.Range("E39").Value = "  +6.57%  "
$ws.Range("D40").Value = "'56.62"
$ws.Range("E40").Value = "  -0.68%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("E42").Value = "  +0.81%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "'3.352.40"
$ws.Range("E43").Value = "  +0.87%  "
$ws.Range("B44").Value = "PEPE"
$ws.Range("C44").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D44").Value = "'0.0₃0734"
$ws.Range("E44").Value = "  +3.99%  "
$ws.Range("E45").Value = "  -4.81%  "
$ws.Range("D46").Value = "'2.89"
$ws.Range("E46").Value = "  -2.58%  "
$ws.Range("D47").Value = "'31.89"
$ws.Range("E47").Value = "  -3.14%  "
$ws.Range("E48").Value = "  -2.32%  "
$ws.Range("E49").Value = "  +0.66%  "
$ws.Range("D50").Value = "'132.80"
$ws.Range("E50").Value = "  -1.58%  "
